$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data: Array of (date serial, value) pairs for rows 2..53
$data = @(
    @(39400, -0.8957516004554691),
    @(39583, 1.910893826230975),
    @(39765, -0.9623800849065276),
    @(39948, -1.734840982136873),
    @(40130, -1.4178810011895),
    @(40310, -1.470147873721189),
    @(40494, 1.938408417700344),
    @(40676, -0.163634158232469),
    @(40862, 0.4),
    @(41044, -0.3),
    @(41228, 0),
    @(41409, 0.009235986179263023),
    @(41592, 0.3863474960573257),
    @(41774, 1.852186157158073),
    @(41957, -0.1726927221574073),
    @(42137, 1.050807574684342),
    @(42321, -0.2184803162966205),
    @(42503, 0.8513884674671885),
    @(42689, 0.4819278240608753),
    @(42867, 0.6629265129002277),
    @(43053, 1.164700738417963),
    @(43145, 0.3865811319001295),
    @(43235, 0.3857269132374052),
    @(43326, 0.4445645077595088),
    @(43418, 0.7),
    @(43510, 0.764905301728362),
    @(43600, 0.2),
    @(43691, -0.1011524282971408),
    @(43783, 0.3399426662647187),
    @(43875, -0.4401201376428787),
    @(43966, -1.5),
    @(44068, -7.231044133207007),
    @(44159, 4.665333711727925),
    @(44251, -0.2672256391354182),
    @(44341, -1.26956845314902),
    @(44432, 2.350117300589673),
    @(44525, 2.031122027629067),
    @(44617, -0.5497577098212645),
    @(44706, 1.775201599534199),
    @(44798, 0.7199425862606432),
    @(44890, 0.5016525134024334),
    @(44981, -0.5890802842202163),
    @(45071, -1.005803602395417),
    @(45163, 0.5789791938159112),
    @(45254, -0.3682387698544858),
    @(45345, -0.2885446736570572),
    @(45436, 0),
    @(45534, 0.009541324901121584),
    @(45618, 1.028868575076984),
    @(45713, 1.007059757688495),
    @(45800, -0.5319104727235242),
    @(45891, 0.4200898674779694)
)

$rowCount = $data.Count
$lastRow = 1 + $rowCount   # row 1 is header

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}

# Remove any leftover rows below the new last data row (previously up to row 73)
$oldLastRow = 73
if ($oldLastRow -gt $lastRow) {
    $deleteRows = $ws.Range($ws.Cells.Item($lastRow + 1, 1), $ws.Cells.Item($oldLastRow, 2)).EntireRow
    $deleteRows.Delete()
}
